$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders ("datetimeFigureOut" fields): 13.08.2018 -> 25.08.2018
#    These live on the slide master, every custom layout and the notes
#    master. Find the shape whose text equals the old date and update it.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "13.08.2018") {
                $tr.Text = "25.08.2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

Update-DatePlaceholder $p.NotesMaster.Shapes

# ---------------------------------------------------------------------
# 2) Slide 1, shape "Rechteck 4" (index 5): "EvoHistorySys" -> "EvoHistorySysGlue"
#    and the leading spacer run shrinks from 20 to 19 spaces.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
$shp5 = $s.Shapes.Item(5)
$tr5 = $shp5.TextFrame.TextRange

$spacer5 = $tr5.Characters(18, 20)
$spacer5.Text = "                   "

$word5 = $tr5.Characters(37, 13)
$word5.Text = "EvoHistorySysGlue"

# ---------------------------------------------------------------------
# 3) Slide 1, shape "Rechteck 21" (index 21): "EvoModelData" -> "EvoModelDataGlue"
#    and the leading spacer run shrinks from 23 to 17 spaces.
# ---------------------------------------------------------------------
$shp21 = $s.Shapes.Item(21)
$tr21 = $shp21.TextFrame.TextRange

$spacer21 = $tr21.Characters(6, 23)
$spacer21.Text = "                 "

$word21 = $tr21.Characters(23, 12)
$word21.Text = "EvoModelDataGlue"
